# Auto-generated edit script: updates Leve profit/price figures per the target diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 167
$ws.Range("H43").Value = 1902.909
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 1694.8572
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 1694.8572
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -1832.8572
$ws.Range("H58").Value = 5081.6523
$ws.Range("I58").Value = 1954.9286
$ws.Range("K58").Value = 5864.7858
$ws.Range("M58").Value = -5714.7858
$ws.Range("H74").Value = 11050.9
$ws.Range("I74").Value = 9501
$ws.Range("K74").Value = 9501
$ws.Range("M74").Value = -8565
$ws.Range("H77").Value = 11050.9
$ws.Range("I77").Value = 9501
$ws.Range("K77").Value = 47505
$ws.Range("M77").Value = -42825
$ws.Range("H86").Value = 2108398.5
$ws.Range("I86").Value = 2800.889
$ws.Range("J86").Value = 3292797.2
$ws.Range("K86").Value = 2800.889
$ws.Range("L86").Value = 3292797.2
$ws.Range("M86").Value = -1677.889
$ws.Range("N86").Value = -3295043.2
$ws.Range("H87").Value = 93500
$ws.Range("J87").Value = 93500
$ws.Range("L87").Value = 93500
$ws.Range("N87").Value = -95996
$ws.Range("H89").Value = 2108398.5
$ws.Range("I89").Value = 2800.889
$ws.Range("J89").Value = 3292797.2
$ws.Range("K89").Value = 14004.445
$ws.Range("L89").Value = 16463986
$ws.Range("M89").Value = -8388.445
$ws.Range("N89").Value = -16475218
$ws.Range("H90").Value = 93500
$ws.Range("J90").Value = 93500
$ws.Range("L90").Value = 280500
$ws.Range("N90").Value = -292980
$ws.Range("H112").Value = 2824.1667
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1848.53
$ws.Range("I32").Value = 1848.53
$ws.Range("K32").Value = 1848.53
$ws.Range("M32").Value = -1561.53
$ws.Range("H88").Value = 2296.9473
$ws.Range("I88").Value = 3832.75
$ws.Range("J88").Value = 1180
$ws.Range("K88").Value = 3832.75
$ws.Range("L88").Value = 1180
$ws.Range("M88").Value = -3426.75
$ws.Range("N88").Value = -1992
$ws.Range("H91").Value = 2296.9473
$ws.Range("I91").Value = 3832.75
$ws.Range("J91").Value = 1180
$ws.Range("K91").Value = 3832.75
$ws.Range("L91").Value = 1180
$ws.Range("M91").Value = -2428.75
$ws.Range("N91").Value = -3988
$ws.Range("H102").Value = 1948.6154
$ws.Range("I102").Value = 1569.3334
$ws.Range("K102").Value = 1569.3334
$ws.Range("M102").Value = 52.66660000000002
$ws.Range("H132").Value = 4244.407
$ws.Range("I132").Value = 4374
$ws.Range("K132").Value = 13122
$ws.Range("M132").Value = -10592

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 710524.0600000001
$ws.Range("J86").Value = 2589.111
$ws.Range("L86").Value = 2589.111
$ws.Range("N86").Value = -4835.111
$ws.Range("H89").Value = 710524.0600000001
$ws.Range("J89").Value = 2589.111
$ws.Range("L89").Value = 12945.555
$ws.Range("N89").Value = -24177.555

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 74013.5
$ws.Range("I52").Value = 75845
$ws.Range("J52").Value = 73555.625
$ws.Range("K52").Value = 75845
$ws.Range("L52").Value = 73555.625
$ws.Range("M52").Value = -75551
$ws.Range("N52").Value = -74143.625
$ws.Range("H99").Value = 4222.75
$ws.Range("I99").Value = 3756.8
$ws.Range("K99").Value = 3756.8
$ws.Range("M99").Value = -2258.8
$ws.Range("H107").Value = 552.36664
$ws.Range("I107").Value = 491.45834
$ws.Range("J107").Value = 796
$ws.Range("K107").Value = 491.45834
$ws.Range("L107").Value = 796
$ws.Range("M107").Value = 1428.54166
$ws.Range("N107").Value = -4636
$ws.Range("H126").Value = 4222.75
$ws.Range("I126").Value = 3756.8
$ws.Range("K126").Value = 11270.4
$ws.Range("M126").Value = -8800.400000000001
$ws.Range("H135").Value = 49582.668
$ws.Range("J135").Value = 49582.668
$ws.Range("L135").Value = 49582.668
$ws.Range("N135").Value = -59722.668
$ws.Range("H137").Value = 58747.5
$ws.Range("J137").Value = 64282.855
$ws.Range("L137").Value = 64282.855
$ws.Range("N137").Value = -74482.85500000001
$ws.Range("H140").Value = 49999.8
$ws.Range("J140").Value = 49999.8
$ws.Range("L140").Value = 49999.8
$ws.Range("N140").Value = -60359.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3000.62
$ws.Range("I131").Value = 1864.1
$ws.Range("J131").Value = 3284.75
$ws.Range("K131").Value = 5592.299999999999
$ws.Range("L131").Value = 9854.25
$ws.Range("M131").Value = -552.2999999999993
$ws.Range("N131").Value = -19934.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3239.7778
$ws.Range("J126").Value = 3879.8572
$ws.Range("L126").Value = 11639.5716
$ws.Range("N126").Value = -16579.5716

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8895.75
$ws.Range("J7").Value = 7798.75
$ws.Range("L7").Value = 7798.75
$ws.Range("N7").Value = -8022.75
$ws.Range("H22").Value = 1351.5333
$ws.Range("I22").Value = 1662.6666
$ws.Range("K22").Value = 1662.6666
$ws.Range("M22").Value = -1367.6666
$ws.Range("H27").Value = 1351.5333
$ws.Range("I27").Value = 1662.6666
$ws.Range("K27").Value = 1662.6666
$ws.Range("M27").Value = -1555.6666
$ws.Range("H68").Value = 2587.7144
$ws.Range("J68").Value = 2774.5
$ws.Range("L68").Value = 2774.5
$ws.Range("N68").Value = -4272.5
$ws.Range("H71").Value = 2587.7144
$ws.Range("J71").Value = 2774.5
$ws.Range("L71").Value = 13872.5
$ws.Range("N71").Value = -21360.5
$ws.Range("H82").Value = 4752.5
$ws.Range("I82").Value = 4755
$ws.Range("J82").Value = 4750
$ws.Range("K82").Value = 4755
$ws.Range("L82").Value = 4750
$ws.Range("M82").Value = -4394
$ws.Range("N82").Value = -5472
$ws.Range("H85").Value = 4752.5
$ws.Range("I85").Value = 4755
$ws.Range("J85").Value = 4750
$ws.Range("K85").Value = 4755
$ws.Range("L85").Value = 4750
$ws.Range("M85").Value = -3507
$ws.Range("N85").Value = -7246
$ws.Range("H126").Value = 8895.75
$ws.Range("J126").Value = 7798.75
$ws.Range("L126").Value = 23396.25
$ws.Range("N126").Value = -28336.25
$ws.Range("H132").Value = 12000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 36000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -41060

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -797
$ws.Range("H81").Value = 4423.1
$ws.Range("I81").Value = 2151.3845
$ws.Range("K81").Value = 4302.769
$ws.Range("M81").Value = -3241.769
$ws.Range("H84").Value = 4423.1
$ws.Range("I84").Value = 2151.3845
$ws.Range("K84").Value = 21513.845
$ws.Range("M84").Value = -16209.845
$ws.Range("H123").Value = 60277.5
$ws.Range("J123").Value = 60277.5
$ws.Range("L123").Value = 60277.5
$ws.Range("N123").Value = -70077.5
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 54722
$ws.Range("J132").Value = 173717.5
$ws.Range("L132").Value = 521152.5
$ws.Range("N132").Value = -526212.5

Write-Host "Applied all Jenova_Profits updates"